$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (pushes current "date" column to G)
$ws.Range("F1").EntireColumn.Insert()

# Set header
$ws.Range("F1").Value = "population"

# Fill population value 354317 for all existing data rows (2 through 361)
$ws.Range("F2:F361").Value = 354317

# Update last two existing rows with new data values
$ws.Range("B360").Value = 18378
$ws.Range("C360").Value = 216
$ws.Range("D360").Value = 7
$ws.Range("E360").Value = 139

$ws.Range("B361").Value = 18465
$ws.Range("C361").Value = 218
$ws.Range("D361").Value = 2
$ws.Range("E361").Value = 87

# Add new rows 362-364 (copy formatting from the last existing date cell first)
$ws.Range("A361").Copy()
$ws.Range("A362:A364").PasteSpecial(-4122)

$ws.Range("A362").Value = "18/mar"
$ws.Range("B362").Value = 18538
$ws.Range("C362").Value = 221
$ws.Range("D362").Value = 3
$ws.Range("E362").Value = 73
$ws.Range("F362").Value = 354317
$ws.Range("G362").Value = "18/mar"

$ws.Range("A363").Value = "19/mar"
$ws.Range("B363").Value = 18650
$ws.Range("C363").Value = 222
$ws.Range("D363").Value = 1
$ws.Range("E363").Value = 112
$ws.Range("F363").Value = 354317
$ws.Range("G363").Value = "19/mar"

$ws.Range("A364").Value = "20/mar"
$ws.Range("B364").Value = 18650
$ws.Range("C364").Value = 222
$ws.Range("D364").Value = 0
$ws.Range("E364").Value = 0
$ws.Range("F364").Value = 354317
$ws.Range("G364").Value = "20/mar"
